$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Add a new worksheet "2022-Q3" right after "总计" (pushes the other
#    quarter sheets down by one position).
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)

$newSheet = $wb.Worksheets.Add($null, $totalSheet)
$newSheet.Name = "2022-Q3"

# ---------------------------------------------------------------------------
# 2) Populate the new "2022-Q3" sheet with the fund-holding data. Columns
#    B, D, E, F, G are stored as plain text in this workbook's convention
#    (fund codes keep leading zeros, percentages/amounts are text) - force
#    Text format while entering them so values aren't coerced to numbers.
# ---------------------------------------------------------------------------
$newSheet.Range("B2:G12").NumberFormat = "@"

$newSheet.Cells.Item(1,2).Value = "基金代码"
$newSheet.Cells.Item(1,3).Value = "基金名称"
$newSheet.Cells.Item(1,4).Value = "基金规模"
$newSheet.Cells.Item(1,5).Value = "股票总仓位"
$newSheet.Cells.Item(1,6).Value = "仓位占比"
$newSheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1,8).Value = "仓位排名"

$rows = @(
  @("003293", "易方达科瑞灵活配置混合", "30.41", "82.90", "3.31", "1.0066", 8),
  @("006533", "易方达科融混合", "24.01", "84.99", "2.97", "0.7131", 10),
  @("110012", "易方达科汇灵活配置混合", "15.60", "74.55", "3.09", "0.4820", 8),
  @("010389", "易方达科益混合A", "6.40", "87.94", "3.58", "0.2291", 7),
  @("011649", "易方达逆向投资混合A", "5.22", "83.77", "3.12", "0.1629", 8),
  @("011650", "易方达逆向投资混合C", "2.59", "83.77", "3.12", "0.0808", 8),
  @("013603", "易方达均衡优选一年持有混合A", "2.48", "49.36", "1.93", "0.0479", 6),
  @("010390", "易方达科益混合C", "0.93", "87.94", "3.58", "0.0333", 7),
  @("013604", "易方达均衡优选一年持有混合C", "0.35", "49.36", "1.93", "0.0068", 6),
  @("005005", "中金金泽量化精选混合A", "0.12", "90.59", "4.09", "0.0049", 5),
  @("005006", "中金金泽量化精选混合C", "0.04", "90.59", "4.09", "0.0016", 5)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $newSheet.Cells.Item($r, 1).Value = $i
    $newSheet.Cells.Item($r, 2).Value = $row[0]
    $newSheet.Cells.Item($r, 3).Value = $row[1]
    $newSheet.Cells.Item($r, 4).Value = $row[2]
    $newSheet.Cells.Item($r, 5).Value = $row[3]
    $newSheet.Cells.Item($r, 6).Value = $row[4]
    $newSheet.Cells.Item($r, 7).Value = $row[5]
    $newSheet.Cells.Item($r, 8).Value = $row[6]
}

# Drop the temporary Text number-format again (kept the stored text values,
# but avoids leaving a non-stock cell style behind).
$newSheet.Range("B2:G12").ClearFormats()

# Borrow header-row / first-column formatting from the existing "2022-Q2"
# sheet (same layout) so the new sheet matches the established look. Look
# this sheet reference up *after* the insert (sheet references resolve by
# position, and inserting a sheet shifts what an old reference points at).
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$q2Sheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$q2Sheet.Range("A2").Copy()
$newSheet.Range("A2:A12").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3) Insert a new leading row in "总计" for 2022-Q3 totals, pushing the
#    older quarters down.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("总计")

# Extend the styled A-column down to the new row 7 (matches rows 2-6).
$ws1.Range("A6").Copy()
$ws1.Range("A7").PasteSpecial(-4122)

$totals = @(
  @("2022-Q3", 11, 2.77),
  @("2022-Q2", 3, 0.59),
  @("2022-Q1", 12, 4.11),
  @("2021-Q4", 7, 2.01),
  @("2021-Q3", 6, 1.67),
  @("2021-Q2", 6, 2.15)
)

for ($i = 0; $i -lt $totals.Count; $i++) {
    $r = $i + 2
    $t = $totals[$i]
    $ws1.Cells.Item($r, 1).Value = $i
    $ws1.Cells.Item($r, 2).Value = $t[0]
    $ws1.Cells.Item($r, 3).Value = $t[1]
    $ws1.Cells.Item($r, 4).Value = $t[2]
}
